$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 778 (old rows 778-806 shift down to 781-809)
$ws.Range("A778:R780").EntireRow.Insert()

# Row 778
$ws.Cells.Item(778,1).Value = 7
$ws.Cells.Item(778,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(778,3).Value = "Ñuble"
$ws.Cells.Item(778,4).Value = 44939
$ws.Cells.Item(778,5).Value = 16
$ws.Cells.Item(778,6).Value = 100112004
$ws.Cells.Item(778,7).Value = "Cebolla"
$ws.Cells.Item(778,8).Value = "Sin especificar"
$ws.Cells.Item(778,9).Value = "1a (cosecha)"
$ws.Cells.Item(778,10).Value = 120
$ws.Cells.Item(778,11).Value = 9000
$ws.Cells.Item(778,12).Value = 9500
$ws.Cells.Item(778,13).Value = 9250
$ws.Cells.Item(778,14).Value = "$/malla 18 kilos"
$ws.Cells.Item(778,15).Value = "Región de O'Higgins"
$ws.Cells.Item(778,16).Value = 514
$ws.Cells.Item(778,17).Value = 18
$ws.Cells.Item(778,18).Value = "Hortaliza"

# Row 779
$ws.Cells.Item(779,1).Value = 7
$ws.Cells.Item(779,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(779,3).Value = "Ñuble"
$ws.Cells.Item(779,4).Value = 44939
$ws.Cells.Item(779,5).Value = 16
$ws.Cells.Item(779,6).Value = 100112004
$ws.Cells.Item(779,7).Value = "Cebolla"
$ws.Cells.Item(779,8).Value = "Sin especificar"
$ws.Cells.Item(779,9).Value = "1a nueva(o)"
$ws.Cells.Item(779,10).Value = 30000
$ws.Cells.Item(779,11).Value = 2100
$ws.Cells.Item(779,12).Value = 2200
$ws.Cells.Item(779,13).Value = 2150
$ws.Cells.Item(779,14).Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Cells.Item(779,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(779,16).Value = 215
$ws.Cells.Item(779,17).Value = 10
$ws.Cells.Item(779,18).Value = "Hortaliza"

# Row 780
$ws.Cells.Item(780,1).Value = 7
$ws.Cells.Item(780,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(780,3).Value = "Ñuble"
$ws.Cells.Item(780,4).Value = 44939
$ws.Cells.Item(780,5).Value = 16
$ws.Cells.Item(780,6).Value = 100112004
$ws.Cells.Item(780,7).Value = "Cebolla"
$ws.Cells.Item(780,8).Value = "Sin especificar"
$ws.Cells.Item(780,9).Value = "2a nueva(o)"
$ws.Cells.Item(780,10).Value = 40000
$ws.Cells.Item(780,11).Value = 1800
$ws.Cells.Item(780,12).Value = 1900
$ws.Cells.Item(780,13).Value = 1850
$ws.Cells.Item(780,14).Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Cells.Item(780,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(780,16).Value = 185
$ws.Cells.Item(780,17).Value = 10
$ws.Cells.Item(780,18).Value = "Hortaliza"
